$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'292.95"
$ws.Range("E2").Value = "'2.18%"
$ws.Range("D3").Value = "'30.60"
$ws.Range("E3").Value = "'7.00%"
$ws.Range("D4").Value = "'5.146"
$ws.Range("E4").Value = "'1.61%"
$ws.Range("D5").Value = "'0.07122"
$ws.Range("E5").Value = "'6.97%"
$ws.Range("D6").Value = "'7.557"
$ws.Range("E6").Value = "'2.30%"
$ws.Range("D7").Value = "'3.628"
$ws.Range("E7").Value = "'6.55%"
$ws.Range("D8").Value = "'1.406"
$ws.Range("E8").Value = "'2.60%"
$ws.Range("D9").Value = "'0.9177"
$ws.Range("E9").Value = "'-1.84%"
$ws.Range("D10").Value = "'0.1628"
$ws.Range("E10").Value = "'3.20%"
$ws.Range("D11").Value = "'0.07677"
$ws.Range("E11").Value = "'16.18%"
$ws.Range("D12").Value = "'0.07778"
$ws.Range("E12").Value = "'2.88%"
$ws.Range("D13").Value = "'0.02945"
$ws.Range("E13").Value = "'0.32%"
$ws.Range("D14").Value = "'0.09009"
$ws.Range("E14").Value = "'0.16%"
$ws.Range("D15").Value = "'0.001594"
$ws.Range("E15").Value = "'0.49%"
$ws.Range("D16").Value = "'0.0006559"
$ws.Range("E16").Value = "'1.37%"
$ws.Range("D17").Value = "'0.006484"
$ws.Range("E17").Value = "'3.62%"
$ws.Range("D18").Value = "'3.482"
$ws.Range("E18").Value = "'1.04%"
$ws.Range("D19").Value = "'2.236"
$ws.Range("E19").Value = "'-0.63%"
$ws.Range("D20").Value = "'0.3277"
$ws.Range("E20").Value = "'1.85%"
$ws.Range("E21").Value = "'5.38%"
$ws.Range("D22").Value = "'3.850"
$ws.Range("E22").Value = "'-5.93%"
$ws.Range("E23").Value = "'3.07%"
$ws.Range("D24").Value = "'0.04537"
$ws.Range("E24").Value = "'0.80%"
$ws.Range("D25").Value = "'0.001211"
$ws.Range("E25").Value = "'2.60%"
$ws.Range("D26").Value = "'0.004245"
$ws.Range("E26").Value = "'2.49%"
$ws.Range("D27").Value = "'0.0001170"
$ws.Range("E27").Value = "'-6.29%"
$ws.Range("D28").Value = "'0.0001689"
$ws.Range("E28").Value = "'4.48%"
$ws.Range("D40").Value = "'0.04412"
$ws.Range("E40").Value = "'4.93%"
$ws.Range("D41").Value = "'0.007028"
$ws.Range("E41").Value = "'4.66%"
$ws.Range("D42").Value = "'0.1276"
$ws.Range("E42").Value = "'2.20%"
$ws.Range("D43").Value = "'0.002210"
$ws.Range("E43").Value = "'9.54%"
$ws.Range("D44").Value = "'0.01328"
$ws.Range("E44").Value = "'9.16%"
$ws.Range("D45").Value = "'0.00005869"
$ws.Range("E45").Value = "'4.37%"
$ws.Range("D47").Value = "'0.01299"
$ws.Range("E47").Value = "'-0.50%"
